$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "70.855.63"
Set-TextValue $ws.Range("E2") "  +3.12%  "

Set-TextValue $ws.Range("D3") "3.571.49"
Set-TextValue $ws.Range("E3") "  +2.23%  "

Set-TextValue $ws.Range("E4") "  +0.09%  "

Set-TextValue $ws.Range("D5") "582.28"
Set-TextValue $ws.Range("E5") "  +2.21%  "

Set-TextValue $ws.Range("D6") "186.59"
Set-TextValue $ws.Range("E6") "  +2.19%  "

Set-TextValue $ws.Range("E7") "  +2.27%  "

Set-TextValue $ws.Range("D8") "3.559.60"
Set-TextValue $ws.Range("E8") "  +2.10%  "

Set-TextValue $ws.Range("E10") "  +22.48%  "

Set-TextValue $ws.Range("D11") "0.651"
Set-TextValue $ws.Range("E11") "  +1.08%  "

Set-TextValue $ws.Range("D12") "54.75"
Set-TextValue $ws.Range("E12") "  +1.67%  "

Set-TextValue $ws.Range("E13") "  +6.77%  "

Set-TextValue $ws.Range("D14") "9.49"
Set-TextValue $ws.Range("E14") "  +0.93%  "

Set-TextValue $ws.Range("D15") "4.141.05"

Set-TextValue $ws.Range("D16") "70.919.93"
Set-TextValue $ws.Range("E16") "  +3.34%  "

Set-TextValue $ws.Range("B17") "WrappedEther"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "3.592.29"
Set-TextValue $ws.Range("E17") "  +3.44%  "

Set-TextValue $ws.Range("B18") "Chainlink"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D18") "19.25"
Set-TextValue $ws.Range("E18") "  +0.19%  "

Set-TextValue $ws.Range("B19") "Uniswap"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D19") "12.84"
Set-TextValue $ws.Range("E19") "  +5.04%  "

Set-TextValue $ws.Range("D20") "574.49"
Set-TextValue $ws.Range("E20") "  +6.33%  "

Set-TextValue $ws.Range("E21") "  +0.64%  "

Set-TextValue $ws.Range("E22") "  -0.38%  "

Set-TextValue $ws.Range("D23") "17.63"
Set-TextValue $ws.Range("E23") "  -6.92%  "

Set-TextValue $ws.Range("D24") "4.57"
Set-TextValue $ws.Range("E24") "  +4.61%  "

Set-TextValue $ws.Range("D25") "4.91"
Set-TextValue $ws.Range("E25") "  -1.94%  "

Set-TextValue $ws.Range("D26") "93.95"
Set-TextValue $ws.Range("E26") "  +0.31%  "

Set-TextValue $ws.Range("D27") "11.25"
Set-TextValue $ws.Range("E27") "  +4.71%  "

Set-TextValue $ws.Range("E28") "  +2.16%  "

Set-TextValue $ws.Range("D29") "9.22"
Set-TextValue $ws.Range("E29") "  +1.27%  "

Set-TextValue $ws.Range("D30") "32.46"
Set-TextValue $ws.Range("E30") "  +3.13%  "

Set-TextValue $ws.Range("D31") "7.22"
Set-TextValue $ws.Range("E31") "  +0.67%  "

Set-TextValue $ws.Range("D32") "12.31"
Set-TextValue $ws.Range("E32") "  -1.28%  "

Set-TextValue $ws.Range("D33") "0.117"
Set-TextValue $ws.Range("E33") "  +3.23%  "

Set-TextValue $ws.Range("D34") "63.14"
Set-TextValue $ws.Range("E34") "  -2.25%  "

Set-TextValue $ws.Range("D35") "3.42"
Set-TextValue $ws.Range("E35") "  +15.27%  "

Set-TextValue $ws.Range("D36") "3.61"
Set-TextValue $ws.Range("E36") "  +16.52%  "

Set-TextValue $ws.Range("D37") "545.19"
Set-TextValue $ws.Range("E37") "  -3.57%  "

Set-TextValue $ws.Range("D38") "0.416"
Set-TextValue $ws.Range("E38") "  +5.68%  "

Set-TextValue $ws.Range("D39") "38.26"
Set-TextValue $ws.Range("E39") "  +1.57%  "

Set-TextValue $ws.Range("D40") "0.0₃0807"
Set-TextValue $ws.Range("E40") "  +5.70%  "

Set-TextValue $ws.Range("D41") "1.00"
Set-TextValue $ws.Range("E41") "  +0.06%  "

Set-TextValue $ws.Range("D42") "3.578.22"
Set-TextValue $ws.Range("E42") "  +10.03%  "

Set-TextValue $ws.Range("E43") "  +5.38%  "

Set-TextValue $ws.Range("E44") "  +3.75%  "

Set-TextValue $ws.Range("D45") "0.0467"
Set-TextValue $ws.Range("E45") "  +7.14%  "

Set-TextValue $ws.Range("E46") "  +0.28%  "

Set-TextValue $ws.Range("D47") "2.93"
Set-TextValue $ws.Range("E47") "  -0.98%  "

Set-TextValue $ws.Range("D48") "9.36"
Set-TextValue $ws.Range("E48") "  +5.10%  "

Set-TextValue $ws.Range("E49") "  +2.85%  "

Set-TextValue $ws.Range("D50") "1.49"
Set-TextValue $ws.Range("E50") "  +14.63%  "

Set-TextValue $ws.Range("E51") "  +0.13%  "

Write-Output "Applied all changes"